$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Mandril row (24-25): predator changes from "Pohon"/"t" to "Beruang"/"b"
$ws.Range("E24").Value = "Beruang"
$ws.Range("E25").Value = "b"

# Clear Beruang row (26-27) predator cells - no longer self-referential
$ws.Range("E26").ClearContents()
$ws.Range("E27").ClearContents()

# Hyena row (30-31): add predator "harimau"/"h"
$ws.Range("D30").Value = "harimau"
$ws.Range("D31").Value = "h"

# Restore the last active selection
$ws.Range("G15").Select()
